$wb = $excel.ActiveWorkbook

# 展览 (sheet1)
$ws = $wb.Worksheets.Item("展览")
$ws.Range("F2").Value = 539
$ws.Range("F3").Value = 930
$ws.Range("F4").Value = 46
$ws.Range("F5").Value = 228
$ws.Range("F7").Value = 1162
$ws.Range("F8").Value = 923
$ws.Range("F10").Value = 725
$ws.Range("F12").Value = 1460
$ws.Range("F13").Value = 60
$ws.Range("F15").Value = 1623
$ws.Range("F16").Value = 19
$ws.Range("F17").Value = 619
$ws.Range("F21").Value = 1085
$ws.Range("F22").Value = 1514
$ws.Range("F24").Value = 627
$ws.Range("F25").Value = 499
$ws.Range("F28").Value = 1023
$ws.Range("F29").Value = 1152
$ws.Range("F30").Value = 310
$ws.Range("F31").Value = 2431
$ws.Range("F33").Value = 1375
$ws.Range("F36").Value = 3973

# 演出 (sheet2)
$ws = $wb.Worksheets.Item("演出")
$ws.Range("F6").Value = 178
$ws.Range("F8").Value = 9
$ws.Range("F9").Value = 7
$ws.Range("F13").Value = 345
$ws.Range("F16").Value = 18
$ws.Range("F17").Value = 27
$ws.Range("F18").Value = 2
$ws.Range("F19").Value = 44
$ws.Range("F21").Value = 259
$ws.Range("F22").Value = 261
$ws.Range("F25").Value = 237
$ws.Range("F29").Value = 1715
$ws.Range("F34").Value = 16

# 本地生活 (sheet3)
$ws = $wb.Worksheets.Item("本地生活")
$ws.Range("F2").Value = 96
$ws.Range("F4").Value = 1272
$ws.Range("F5").Value = 1669
$ws.Range("F6").Value = 452
$ws.Range("F7").Value = 1012

# 全部类型 (sheet4)
$ws = $wb.Worksheets.Item("全部类型")
$ws.Range("F2").Value = 96
$ws.Range("F4").Value = 1272
$ws.Range("F5").Value = 1669
$ws.Range("F6").Value = 452
$ws.Range("F7").Value = 1012
$ws.Range("F8").Value = 539
$ws.Range("F9").Value = 930
$ws.Range("F10").Value = 46
$ws.Range("F11").Value = 228
$ws.Range("F13").Value = 1162
$ws.Range("F14").Value = 923
$ws.Range("F18").Value = 725
$ws.Range("F19").Value = 178
$ws.Range("F20").Value = 178
$ws.Range("F21").Value = 9
$ws.Range("F23").Value = 1460
$ws.Range("F24").Value = 60
$ws.Range("F26").Value = 1624
$ws.Range("F27").Value = 19
$ws.Range("F28").Value = 619
$ws.Range("F29").Value = 345
$ws.Range("F31").Value = 1085
$ws.Range("F32").Value = 1514
$ws.Range("F34").Value = 627
$ws.Range("F35").Value = 499
$ws.Range("F38").Value = 44
$ws.Range("F39").Value = 261
$ws.Range("F40").Value = 1023
$ws.Range("F41").Value = 1152
$ws.Range("F42").Value = 310
$ws.Range("F43").Value = 2431
$ws.Range("F44").Value = 237
$ws.Range("F46").Value = 1715
$ws.Range("F47").Value = 1715
$ws.Range("F48").Value = 1375
$ws.Range("F51").Value = 3973
